# The "От <tab> г." paragraph (split across 5 runs: "О","т"," ",<tab>,"г.")
# becomes a single run "от____________________________20____г." at 11pt
# (half-points 24 -> 22) instead of 12pt (sz 24).
$d = $word.ActiveDocument

$rng = $d.Content
$searchText = "От " + [char]9 + "г."
$replacementText = "от____________________________20____" + [char]0x0433 + "."

$rng.Find.ClearFormatting()
$found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replacementText, 2)

if ($found) {
    $rng.Font.Size = 11
} else {
    # Fallback: locate the paragraph by its partial text and rebuild it manually
    # (excluding the trailing paragraph mark from the replaced range).
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "От*г.*") {
            $r = $p.Range
            $r.End = $r.End - 1
            $r.Text = $replacementText
            $r.Font.Size = 11
            break
        }
    }
}
